$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value2 = 928561.1
$ws.Range("J17").Value2 = 1092072.1
$ws.Range("L17").Value2 = 3276216.3
$ws.Range("N17").Value2 = -3276552.3
$ws.Range("H40").Value2 = 1384.4445
$ws.Range("I40").Value2 = 1260.1666
$ws.Range("J40").Value2 = 1633
$ws.Range("K40").Value2 = 1260.1666
$ws.Range("L40").Value2 = 1633
$ws.Range("M40").Value2 = -1085.1666
$ws.Range("N40").Value2 = -1983
$ws.Range("H63").Value2 = 0
$ws.Range("J63").Value2 = 0
$ws.Range("L63").ClearContents()
$ws.Range("N63").Value2 = 0
$ws.Range("H64").Value2 = 3808.6
$ws.Range("I64").Value2 = 3200
$ws.Range("K64").Value2 = 3200
$ws.Range("M64").Value2 = -2952
$ws.Range("H66").Value2 = 0
$ws.Range("J66").Value2 = 0
$ws.Range("L66").ClearContents()
$ws.Range("N66").Value2 = 0
$ws.Range("H67").Value2 = 3808.6
$ws.Range("I67").Value2 = 3200
$ws.Range("K67").Value2 = 3200
$ws.Range("M67").Value2 = -2342
$ws.Range("H107").Value2 = 1787.7059
$ws.Range("I107").Value2 = 872.96155
$ws.Range("K107").Value2 = 872.96155
$ws.Range("M107").Value2 = 1047.03845
$ws.Range("H125").Value2 = 2592.8
$ws.Range("I125").Value2 = 2592.8
$ws.Range("K125").Value2 = 23335.2
$ws.Range("M125").Value2 = -20875.2
$ws.Range("H135").Value2 = 16670650
$ws.Range("I135").Value2 = 18522390
$ws.Range("J135").Value2 = 4995
$ws.Range("K135").Value2 = 166701510
$ws.Range("L135").Value2 = 44955
$ws.Range("M135").Value2 = -166698975
$ws.Range("N135").Value2 = -50025
$ws.Range("H138").Value2 = 3361.927
$ws.Range("I138").Value2 = 1332.826
$ws.Range("J138").Value2 = 4001.233
$ws.Range("K138").Value2 = 3998.478
$ws.Range("L138").Value2 = 12003.699
$ws.Range("M138").Value2 = 1141.522
$ws.Range("N138").Value2 = -22283.699
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 47403.523
$ws.Range("J32").Value2 = 67305.31
$ws.Range("L32").Value2 = 67305.31
$ws.Range("N32").Value2 = -67879.31
$ws.Range("H44").Value2 = 49500
$ws.Range("J44").Value2 = 49500
$ws.Range("L44").Value2 = 49500
$ws.Range("N44").Value2 = -50476
$ws.Range("H55").Value2 = 21632.666
$ws.Range("J55").Value2 = 27499.5
$ws.Range("L55").Value2 = 27499.5
$ws.Range("N55").Value2 = -28129.5
$ws.Range("H61").Value2 = 3752.3447
$ws.Range("I61").Value2 = 3636.3572
$ws.Range("K61").Value2 = 3636.3572
$ws.Range("M61").Value2 = -3424.3572
$ws.Range("H97").Value2 = 2212.8333
$ws.Range("I97").Value2 = 2239.4375
$ws.Range("K97").Value2 = 2239.4375
$ws.Range("M97").Value2 = -1743.4375
$ws.Range("H132").Value2 = 3551.6667
$ws.Range("I132").Value2 = 2897.7273
$ws.Range("J132").Value2 = 4579.2856
$ws.Range("K132").Value2 = 8693.1819
$ws.Range("L132").Value2 = 13737.8568
$ws.Range("M132").Value2 = -6163.1819
$ws.Range("N132").Value2 = -18797.8568
$ws.Range("H136").Value2 = 3752.3447
$ws.Range("I136").Value2 = 3636.3572
$ws.Range("K136").Value2 = 10909.0716
$ws.Range("M136").Value2 = -8359.071599999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value2 = 1155.2632
$ws.Range("I7").Value2 = 225
$ws.Range("J7").Value2 = 1264.7059
$ws.Range("K7").Value2 = 225
$ws.Range("L7").Value2 = 1264.7059
$ws.Range("M7").Value2 = -112
$ws.Range("N7").Value2 = -1490.7059
$ws.Range("H20").Value2 = 4598.5
$ws.Range("I20").Value2 = 4423.5264
$ws.Range("J20").Value2 = 4967.8887
$ws.Range("K20").Value2 = 4423.5264
$ws.Range("L20").Value2 = 4967.8887
$ws.Range("M20").Value2 = -4176.5264
$ws.Range("N20").Value2 = -5461.8887
$ws.Range("H94").Value2 = 1637.5714
$ws.Range("I94").Value2 = 1575.5
$ws.Range("K94").Value2 = 1575.5
$ws.Range("M94").Value2 = -1124.5
$ws.Range("H107").Value2 = 1640.2727
$ws.Range("I107").Value2 = 1741.6666
$ws.Range("K107").Value2 = 1741.6666
$ws.Range("M107").Value2 = 178.3334
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value2 = 2844.25
$ws.Range("I58").Value2 = 2267.5293
$ws.Range("K58").Value2 = 2267.5293
$ws.Range("M58").Value2 = -2064.5293
$ws.Range("H63").Value2 = 0
$ws.Range("J63").Value2 = 0
$ws.Range("L63").ClearContents()
$ws.Range("N63").Value2 = 0
$ws.Range("H64").Value2 = 13000
$ws.Range("J64").Value2 = 13000
$ws.Range("L64").Value2 = 13000
$ws.Range("N64").Value2 = -13496
$ws.Range("H66").Value2 = 0
$ws.Range("J66").Value2 = 0
$ws.Range("L66").ClearContents()
$ws.Range("N66").Value2 = 0
$ws.Range("H67").Value2 = 13000
$ws.Range("J67").Value2 = 13000
$ws.Range("L67").Value2 = 13000
$ws.Range("N67").Value2 = -14716
$ws.Range("H69").Value2 = 19769.8
$ws.Range("I69").Value2 = 6250
$ws.Range("K69").Value2 = 6250
$ws.Range("M69").Value2 = -5501
$ws.Range("H72").Value2 = 19769.8
$ws.Range("I72").Value2 = 6250
$ws.Range("K72").Value2 = 18750
$ws.Range("M72").Value2 = -15006
$ws.Range("H132").Value2 = 1919.35
$ws.Range("I132").Value2 = 1919.35
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 5758.049999999999
$ws.Range("L132").Value2 = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value2 = -3228.049999999999
$ws.Range("H136").Value2 = 2844.25
$ws.Range("I136").Value2 = 2267.5293
$ws.Range("K136").Value2 = 6802.5879
$ws.Range("M136").Value2 = -4252.5879
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value2 = 2975.75
$ws.Range("J55").Value2 = 3440.5
$ws.Range("L55").Value2 = 10321.5
$ws.Range("N55").Value2 = -10675.5
$ws.Range("I88").Value2 = 5608
$ws.Range("J88").Value2 = 100000000
$ws.Range("K88").Value2 = 16824
$ws.Range("L88").Value2 = 300000000
$ws.Range("M88").Value2 = -16396
$ws.Range("N88").Value2 = -300000856
$ws.Range("I91").Value2 = 5608
$ws.Range("J91").Value2 = 100000000
$ws.Range("K91").Value2 = 16824
$ws.Range("L91").Value2 = 300000000
$ws.Range("M91").Value2 = -15342
$ws.Range("N91").Value2 = -300002964
$ws.Range("H129").Value2 = 5824390
$ws.Range("I129").Value2 = 8250428
$ws.Range("J129").Value2 = 1900
$ws.Range("K129").Value2 = 24751284
$ws.Range("L129").Value2 = 5700
$ws.Range("M129").Value2 = -24746284
$ws.Range("N129").Value2 = -15700
$ws.Range("H131").Value2 = 20361.564
$ws.Range("I131").Value2 = 78734.234
$ws.Range("J131").Value2 = 2293.8333
$ws.Range("K131").Value2 = 236202.702
$ws.Range("L131").Value2 = 6881.499899999999
$ws.Range("M131").Value2 = -231162.702
$ws.Range("N131").Value2 = -16961.4999
$ws.Range("H137").Value2 = 6513.684
$ws.Range("J137").Value2 = 7666.6665
$ws.Range("L137").Value2 = 22999.9995
$ws.Range("N137").Value2 = -33199.99950000001
$ws.Range("H141").Value2 = 335354.34
$ws.Range("I141").Value2 = 3015
$ws.Range("J141").Value2 = 1000033
$ws.Range("K141").Value2 = 9045
$ws.Range("L141").Value2 = 3000099
$ws.Range("M141").Value2 = -3865
$ws.Range("N141").Value2 = -3010459
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value2 = 0
$ws.Range("J64").Value2 = 0
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value2 = 0
$ws.Range("H67").Value2 = 0
$ws.Range("J67").Value2 = 0
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value2 = 0
$ws.Range("H69").Value2 = 77200
$ws.Range("J69").Value2 = 77200
$ws.Range("L69").Value2 = 77200
$ws.Range("N69").Value2 = -78698
$ws.Range("H72").Value2 = 77200
$ws.Range("J72").Value2 = 77200
$ws.Range("L72").Value2 = 231600
$ws.Range("N72").Value2 = -239088
$ws.Range("H97").Value2 = 1612.8276
$ws.Range("I97").Value2 = 1576.3334
$ws.Range("K97").Value2 = 1576.3334
$ws.Range("M97").Value2 = -1080.3334
$ws.Range("H135").Value2 = 49567
$ws.Range("J135").Value2 = 49567
$ws.Range("L135").Value2 = 49567
$ws.Range("N135").Value2 = -59707
$ws.Range("H136").Value2 = 35102
$ws.Range("J136").Value2 = 35102
$ws.Range("L136").Value2 = 105306
$ws.Range("N136").Value2 = -110406
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value2 = 6416.6665
$ws.Range("I61").Value2 = 5625
$ws.Range("K61").Value2 = 5625
$ws.Range("M61").Value2 = -5423
$ws.Range("H113").Value2 = 6416.6665
$ws.Range("I113").Value2 = 5625
$ws.Range("K113").Value2 = 5625
$ws.Range("M113").Value2 = -3455
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 872.88464
$ws.Range("I107").Value2 = 718.2105
$ws.Range("J107").Value2 = 1292.7142
$ws.Range("K107").Value2 = 2154.6315
$ws.Range("L107").Value2 = 3878.1426
$ws.Range("M107").Value2 = -234.6315
$ws.Range("N107").Value2 = -7718.142599999999
$ws.Range("H113").Value2 = 1358.6
$ws.Range("I113").Value2 = 1570.8334
$ws.Range("J113").Value2 = 812.8570999999999
$ws.Range("K113").Value2 = 4712.5002
$ws.Range("L113").Value2 = 2438.5713
$ws.Range("M113").Value2 = -2542.5002
$ws.Range("N113").Value2 = -6778.5713
$ws.Range("H132").Value2 = 4395.1904
$ws.Range("I132").Value2 = 4478.8945
$ws.Range("K132").Value2 = 13436.6835
$ws.Range("M132").Value2 = -10906.6835
$ws.Range("H136").Value2 = 2242.875
$ws.Range("I136").Value2 = 1929.3077
$ws.Range("K136").Value2 = 5787.9231
$ws.Range("M136").Value2 = -3237.9231
